$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 141, shifting existing rows 141-145 down to 142-146
$ws.Rows.Item(141).Insert()

# Populate the newly inserted row 141 with the new weekly price record
$ws.Cells.Item(141, 1).Value = 7
$ws.Cells.Item(141, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(141, 3).Value = "Ñuble"
$ws.Cells.Item(141, 4).Value = 44509
$ws.Cells.Item(141, 5).Value = 16
$ws.Cells.Item(141, 6).Value = 100112003
$ws.Cells.Item(141, 7).Value = "Ajo"
$ws.Cells.Item(141, 8).Value = "Chino"
$ws.Cells.Item(141, 9).Value = "Primera"
$ws.Cells.Item(141, 10).Value = 60
$ws.Cells.Item(141, 11).Value = 15000
$ws.Cells.Item(141, 12).Value = 16000
$ws.Cells.Item(141, 13).Value = 15500
$ws.Cells.Item(141, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(141, 15).Value = "China"
$ws.Cells.Item(141, 16).Value = 1550
$ws.Cells.Item(141, 17).Value = 10
$ws.Cells.Item(141, 18).Value = "Hortaliza"
